# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Rows 16-52 hold the "EMEL OBREGON ORTEGA" (CC 9168686) dues periods,
# previously listed oldest->newest (1607..1901, with row 47 interrupted
# by an "ORLANDO DE JESUS ROJAS DONADO" row for period 1902). The sheet
# is rebuilt so the periods run newest->oldest (1907..1607) in rows
# 16-52, and the lone ORLANDO row (period 1902) moves down to row 53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 previously belonged to ORLANDO DE JESUS ROJAS DONADO; restore it
# to the standard EMEL OBREGON ORTEGA worker row (same Tipo/N° Doc/Nombre
# as the rest of the block) before re-numbering the periods below.
$ws.Range("C47").Value = "9168686"
$ws.Range("D47").Value = "EMEL OBREGON ORTEGA"

# New descending period order (1907 -> 1607) for rows 16-52.
$periods = @(
  "1907","1906","1905","1904","1903","1902","1901",
  "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
  "1712","1711","1710","1709","1708","1707","1706","1705","1704","1703","1702","1701",
  "1612","1611","1610","1609","1608","1607"
)

$row = 16
foreach ($p in $periods) {
    $ws.Range("E$row").Value = $p
    $row++
}

# Valor Mora (F) stays 48000 for every period except the newest (1907,
# row 16), which drops to 24000. Salario Basico (G) stays 1200000.
$ws.Range("F16").Value = 24000
$ws.Range("F17:F52").Value = 48000
$ws.Range("G16:G52").Value = 1200000

# Row 53 now carries the ORLANDO DE JESUS ROJAS DONADO record that used
# to sit at row 47 (period 1902, Valor Mora 32000, Salario Basico 828116).
$ws.Range("C53").Value = "9166846"
$ws.Range("D53").Value = "ORLANDO DE JESUS ROJAS DONADO"
$ws.Range("E53").Value = "1902"
$ws.Range("F53").Value = 32000
$ws.Range("G53").Value = 828116
